# Update the Minneapolis GDP data: revise existing GDP figures (2002-2019)
# and append a new observation row for 2020-01-01, matching a refreshed
# FRED data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised GDP values for existing rows (B12:B30)
$updates = @{
    12 = 140695.81899999999
    13 = 144892.53899999999
    14 = 154369.671
    15 = 165416.53
    16 = 176747.829
    17 = 180933.17
    18 = 188527.274
    19 = 190205.06899999999
    20 = 186323.291
    21 = 193513.28700000001
    22 = 202332.22
    23 = 210304.679
    24 = 219793.32699999999
    25 = 232315.465
    26 = 241312.26300000001
    27 = 247968.00899999999
    28 = 256332.87100000001
    29 = 271255.00400000002
    30 = 278263.82500000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# Append the new observation row (row 31): date 2020-01-01 with matching
# date style, and its GDP value with the matching numeric style.
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").Value = 43831

$ws.Range("B30").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("B31").Value = 270282.12800000003

# Update the active selection/cursor to rest on the newly-added row
# (mirrors the workbook being left with the new observation selected).
$ws.Range("A31").Select() | Out-Null
